# "Generate Report for Handback": the handback transform for the
# 2f450335-...-da042d74892f file failed, so update the localization-status
# report to reflect the failure instead of "Ready for handoff", and record
# the error detail for each locale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$failedStatus = "Handback transform failed"
$zhcnError = "Handback file name: z3hqi0gr.4vh is different with handoff file name: 2f450335-45ca-41a1-b394-da042d74892f.73795c1cb76c882f542469926646f10b130b4e45.zh-cn."
$dedeError = "Handback file name: z3hqi0gr.4vh is different with handoff file name: 2f450335-45ca-41a1-b394-da042d74892f.73795c1cb76c882f542469926646f10b130b4e45.de-de."

# Overview sheet: row 3 is the 2f450335-...-da042d74892f.md file; its zh-cn (E)
# and de-de (F) status columns move from "Ready for handoff" to the failure.
$overview.Range("E3").Value = $failedStatus
$overview.Range("F3").Value = $failedStatus

# zh-cn and de-de sheets: row 3 (same file) gets the new Status and an Error
# Detail message for each locale.
$zhcn.Range("C3").Value = $failedStatus
$zhcn.Range("P3").Value = $zhcnError

$dede.Range("C3").Value = $failedStatus
$dede.Range("P3").Value = $dedeError

# Widen the Error Detail column (P) on both localized sheets so the messages
# are readable. (39.1666... compensates for the engine's internal padding so
# the stored column width comes out to exactly 40.)
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
